$ws = $excel.ActiveWorkbook.ActiveSheet

# Row data for rows 2-17 (Sending cluster x Target cluster, 4x4 matrix of ECs/FAPs/MuSCs/Resolving-Mac)
$rows = @(
    @("ECs", "Slit2", "Robo1", "ECs", "3", "1", "0.2160193333333333", "0.648058", "0.02486881244588016", "0.02486881244588016", "3", "1", "0.7788713333333334", "2.336614", "0.02154486428780801", "0.02154486428780801", "0.1682512661791111", "1.514261395612", "0.0005357951891454387", "0.0005357951891454388"),
    @("ECs", "Slit2", "Robo1", "FAPs", "3", "1", "0.2160193333333333", "0.648058", "0.02486881244588016", "0.02486881244588016", "3", "1", "29.811843", "89.435529", "0.8246446930529893", "0.8246446930529893", "6.439934450298", "57.959410052682", "0.0205079342060252", "0.02050793420602521"),
    @("ECs", "Slit2", "Robo1", "MuSCs", "3", "1", "0.2160193333333333", "0.648058", "0.02486881244588016", "0.02486881244588016", "3", "1", "5.508900000000001", "16.5267", "0.152385250035015", "0.152385250035015", "1.1900289054", "10.7102601486", "0.003789640202639342", "0.003789640202639343"),
    @("ECs", "Slit2", "Robo1", "Resolving-Mac", "3", "1", "0.2160193333333333", "0.648058", "0.02486881244588016", "0.02486881244588016", "1", "0.3333333333333333", "0.05152233333333334", "0.154567", "0.001425192624187658", "0.001425192624187658", "0.01112982009844445", "0.100168380886", "3.544284807017464E-05", "3.544284807017464E-05"),
    @("FAPs", "Slit2", "Robo1", "ECs", "3", "1", "8.296819666666666", "24.890459", "0.9551554900377276", "0.9551554900377278", "3", "1", "0.7788713333333334", "2.336614", "0.02154486428780801", "0.02154486428780801", "6.462154996202889", "58.15939496582599", "0.0205786954066176", "0.0205786954066176"),
    @("FAPs", "Slit2", "Robo1", "FAPs", "3", "1", "8.296819666666666", "24.890459", "0.9551554900377276", "0.9551554900377278", "3", "1", "29.811843", "89.435529", "0.8246446930529893", "0.8246446930529893", "247.343485301979", "2226.091367717811", "0.7876639059000394", "0.7876639059000397"),
    @("FAPs", "Slit2", "Robo1", "MuSCs", "3", "1", "8.296819666666666", "24.890459", "0.9551554900377276", "0.9551554900377278", "3", "1", "5.508900000000001", "16.5267", "0.152385250035015", "0.152385250035015", "45.7063498617", "411.3571487553", "0.1455516081717164", "0.1455516081717165"),
    @("FAPs", "Slit2", "Robo1", "Resolving-Mac", "3", "1", "8.296819666666666", "24.890459", "0.9551554900377276", "0.9551554900377278", "1", "0.3333333333333333", "0.05152233333333334", "0.154567", "0.001425192624187658", "0.001425192624187658", "0.4274715084725556", "3.847243576253", "0.001361280559354118", "0.001361280559354118"),
    @("MuSCs", "Slit2", "Robo1", "ECs", "3", "1", "0.173491", "0.520473", "0.01997281943922393", "0.01997281943922393", "3", "1", "0.7788713333333334", "2.336614", "0.02154486428780801", "0.02154486428780801", "0.1351271664913333", "1.216144498422", "0.0004303116842629732", "0.0004303116842629733"),
    @("MuSCs", "Slit2", "Robo1", "FAPs", "3", "1", "0.173491", "0.520473", "0.01997281943922393", "0.01997281943922393", "3", "1", "29.811843", "89.435529", "0.8246446930529893", "0.8246446930529893", "5.172086453912999", "46.54877808521699", "0.01647047955586159", "0.0164704795558616"),
    @("MuSCs", "Slit2", "Robo1", "MuSCs", "3", "1", "0.173491", "0.520473", "0.01997281943922393", "0.01997281943922393", "3", "1", "5.508900000000001", "16.5267", "0.152385250035015", "0.152385250035015", "0.9557445699", "8.6017011291", "0.003043563084150347", "0.003043563084150348"),
    @("MuSCs", "Slit2", "Robo1", "Resolving-Mac", "3", "1", "0.173491", "0.520473", "0.01997281943922393", "0.01997281943922393", "1", "0.3333333333333333", "0.05152233333333334", "0.154567", "0.001425192624187658", "0.001425192624187658", "0.008938661132333332", "0.080447950191", "2.846511494901382E-05", "2.846511494901383E-05"),
    @("Resolving-Mac", "Slit2", "Robo1", "ECs", "1", "0.3333333333333333", "2.5E-05", "7.499999999999999E-05", "2.878077168156263E-06", "2.878077168156264E-06", "3", "1", "0.7788713333333334", "2.336614", "0.02154486428780801", "0.02154486428780801", "1.947178333333333E-05", "0.00017524605", "6.200778199776549E-08", "6.200778199776549E-08"),
    @("Resolving-Mac", "Slit2", "Robo1", "FAPs", "1", "0.3333333333333333", "2.5E-05", "7.499999999999999E-05", "2.878077168156263E-06", "2.878077168156264E-06", "3", "1", "29.811843", "89.435529", "0.8246446930529893", "0.8246446930529893", "0.0007452960749999999", "0.006707664675", "2.373391062917038E-06", "2.373391062917039E-06"),
    @("Resolving-Mac", "Slit2", "Robo1", "MuSCs", "1", "0.3333333333333333", "2.5E-05", "7.499999999999999E-05", "2.878077168156263E-06", "2.878077168156264E-06", "3", "1", "5.508900000000001", "16.5267", "0.152385250035015", "0.152385250035015", "0.0001377225", "0.0012395025", "4.385765088895602E-07", "4.385765088895603E-07"),
    @("Resolving-Mac", "Slit2", "Robo1", "Resolving-Mac", "1", "0.3333333333333333", "2.5E-05", "7.499999999999999E-05", "2.878077168156263E-06", "2.878077168156264E-06", "1", "0.3333333333333333", "0.05152233333333334", "0.154567", "0.001425192624187658", "0.001425192624187658", "1.288058333333333E-06", "1.1592525E-05", "4.101814351899209E-09", "4.10181435189921E-09")
)

$stringCols = @(1,2,3,4)  # A,B,C,D are 1-based column indices that are strings

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowVals = $rows[$i]
    for ($c = 1; $c -le 20; $c++) {
        $val = $rowVals[$c - 1]
        if ($stringCols -contains $c) {
            $ws.Cells.Item($r, $c).Value = $val
        } else {
            $ws.Cells.Item($r, $c).Value = [double]$val
        }
    }
}

Write-Output "Done writing rows"